# Update the RF (column I) values for rows 23-44 on the active sheet,
# reflecting the "Update of 2025 data and RF changes" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I23:I44").Value2 = 9.814259259259259
